$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match column width of the new column V to column U
$ws.Range("V1").ColumnWidth = $ws.Range("U1").ColumnWidth

# Row 1: date header text (use formula trick to avoid auto date-parsing turning it into a serial number)
$ws.Range("V1").Formula = '="2026/01/05"'

# Row 2: index name text
$ws.Range("V2").Value = "上证"

# Rows 3-118: numeric data values
$ws.Range("V3").Value = 63.73
$ws.Range("V4").Value = 4011.45
$ws.Range("V6").Value = 50.5
$ws.Range("V7").Value = 5745.04
$ws.Range("V9").Value = 54.51
$ws.Range("V10").Value = 4703.41
$ws.Range("V12").Value = 58.37
$ws.Range("V13").Value = 7626.72
$ws.Range("V15").Value = 28.09
$ws.Range("V16").Value = 2713.89
$ws.Range("V18").Value = 96.52
$ws.Range("V19").Value = 6858.47
$ws.Range("V21").Value = 71.59999999999999
$ws.Range("V22").Value = 85700.96000000001
$ws.Range("V24").Value = 86.12
$ws.Range("V25").Value = 24539.34
$ws.Range("V27").Value = 70.19
$ws.Range("V28").Value = 51850.54
$ws.Range("V30").Value = 55.58
$ws.Range("V31").Value = 5520.95
$ws.Range("V33").Value = 1.28
$ws.Range("V34").Value = 32269.26
$ws.Range("V36").Value = 30.09
$ws.Range("V37").Value = 3462
$ws.Range("V39").Value = 49.05
$ws.Range("V40").Value = 3272.07
$ws.Range("V42").Value = 15.89
$ws.Range("V43").Value = 7067.29
$ws.Range("V45").Value = 25.21
$ws.Range("V46").Value = 8515.639999999999
$ws.Range("V48").Value = 7.12
$ws.Range("V49").Value = 12740.56
$ws.Range("V51").Value = 28.2
$ws.Range("V52").Value = 12426.73
$ws.Range("V54").Value = 16.84
$ws.Range("V55").Value = 8914.65
$ws.Range("V57").Value = 24.16
$ws.Range("V58").Value = 15403.84
$ws.Range("V60").Value = 30.06
$ws.Range("V61").Value = 15054.07
$ws.Range("V63").Value = 20.02
$ws.Range("V64").Value = 9765.02
$ws.Range("V66").Value = 11.53
$ws.Range("V67").Value = 9661.379999999999
$ws.Range("V69").Value = 22.96
$ws.Range("V70").Value = 3065.3
$ws.Range("V72").Value = 43.99
$ws.Range("V73").Value = 5726.11
$ws.Range("V75").Value = 23.36
$ws.Range("V76").Value = 9037.35
$ws.Range("V78").Value = 17.02
$ws.Range("V79").Value = 1505.88
$ws.Range("V81").Value = 16.97
$ws.Range("V82").Value = 16217.81
$ws.Range("V84").Value = 52.62
$ws.Range("V85").Value = 2865.99
$ws.Range("V87").Value = 58.48
$ws.Range("V88").Value = 2903.95
$ws.Range("V90").Value = 51.27
$ws.Range("V91").Value = 3103.22
$ws.Range("V93").Value = 42.48
$ws.Range("V94").Value = 1970.07
$ws.Range("V96").Value = 25.96
$ws.Range("V97").Value = 10292.69
$ws.Range("V99").Value = 86.15000000000001
$ws.Range("V100").Value = 9572.860000000001
$ws.Range("V102").Value = 58.2
$ws.Range("V103").Value = 14411.79
$ws.Range("V105").Value = 6.36
$ws.Range("V106").Value = 2215.71
$ws.Range("V108").Value = 25.71
$ws.Range("V109").Value = 865.83
$ws.Range("V111").Value = 31.51
$ws.Range("V112").Value = 3058.48
$ws.Range("V114").Value = 20.84
$ws.Range("V115").Value = 3920.55
$ws.Range("V117").Value = 29.02
$ws.Range("V118").Value = 2947.02

# Apply formatting to match column U (font size 12, centered horizontally/vertically); row 2 is bold like U2
$dataRange = $ws.Range("V1:V118")
$dataRange.Font.Size = 12
$dataRange.HorizontalAlignment = -4108
$dataRange.VerticalAlignment = -4108
$ws.Range("V2").Font.Bold = $true

Write-Output "done"